# Update antiSmash UI requirement input to allow either fasta or genbank
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E4").Value = "Input FASTA/Genbank File"
$ws.Range("A4").Value = "File inputFile"

# Column E should widen (best-fit) to accommodate the longer text now stored there
$ws.Columns("E").EntireColumn.ColumnWidth = 63.94

# Restore the cursor/selection position as recorded in the saved file
$ws.Range("A4").Select()
